{"js": "// The last editing session removed the \"Time: 0 Hours\" paragraph that used\n// to sit right under the chapter title. Word automatically tracks the\n// location of the most recent edit with the hidden \"_GoBack\" bookmark, so\n// after removing that paragraph the bookmark also needs to move from its\n// old spot (end of the doc, on the page-break paragraph) to the start of\n// the paragraph that is now right after the title (\"This chapter covers...\").\n\n// 1) Remove the existing \"_GoBack\" bookmark (wherever Word last left it).\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Locate the \"Time: 0 Hours\" paragraph and the paragraph that follows it.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet targetParagraph = null;\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text === \"Time: 0 Hours\") {\n    targetParagraph = paragraph;\n    break;\n  }\n}\n\nif (!targetParagraph) {\n  throw new Error('Could not find the \"Time: 0 Hours\" paragraph.');\n}\n\nconst nextParagraph = targetParagraph.getNext();\nnextParagraph.load(\"text\");\nawait context.sync();\n\n// 3) Delete the whole \"Time: 0 Hours\" paragraph (including its paragraph\n//    mark), which merges it away and leaves \"This chapter covers...\"\n//    directly under the title.\ntargetParagraph.delete();\nawait context.sync();\n\n// 4) Re-insert the \"_GoBack\" bookmark, collapsed at the very start of the\n//    paragraph that now marks the most recent edit location.\nconst startOfNextParagraph = nextParagraph.getRange(\"Start\");\nstartOfNextParagraph.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The last editing session removed the \"Time: 0 Hours\" paragraph that used\n# to sit right under the chapter title. Word automatically tracks the\n# location of the most recent edit with the hidden \"_GoBack\" bookmark, so\n# after removing that paragraph the bookmark also needs to move from its\n# old spot (end of the doc, on the page-break paragraph) to the start of\n# the paragraph that is now right after the title (\"This chapter covers...\").\n\n$d = $word.ActiveDocument\n\n# 1) Remove the existing \"_GoBack\" bookmark (wherever Word last left it).\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Locate the \"Time: 0 Hours\" paragraph and delete it entirely (including\n#    its paragraph mark), which merges it away and leaves \"This chapter\n#    covers...\" directly under the title.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Time: 0 Hours\"\n$found = $find.Execute()\nif (-not $found) {\n    throw \"Could not find the 'Time: 0 Hours' paragraph\"\n}\n$timeParagraph = $find.Parent.Paragraphs(1)\n$timeParagraph.Range.Delete()\n\n# 3) Re-locate the paragraph that is now the most recently edited spot (the\n#    \"This chapter covers...\" paragraph, now promoted right after the\n#    title) and drop a collapsed \"_GoBack\" bookmark at its very start.\n#    NOTE: paragraph/range objects captured before the delete above point\n#    at stale offsets once the document shrinks, so we search again.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Text = \"This chapter covers\"\n$found2 = $find2.Execute()\nif (-not $found2) {\n    throw \"Could not find the 'This chapter covers' paragraph\"\n}\n$nextParagraph = $find2.Parent.Paragraphs(1)\n$bookmarkStart = $nextParagraph.Range.Start\n$bookmarkRange = $d.Range($bookmarkStart, $bookmarkStart)\n$d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n"}
